# Commit message: Change "Bauteil" to "Component" Change "bauteil" to "component"
#
# The workbook "ExcelVorlage.xlsx" has three sheets (IfcWall, IfcWindow,
# IfcDoor). Each sheet uses the German column header "Bauteil" in cell A1.
# Rename that header to the English "Component" on every sheet (case
# insensitive, in case the original text was lower-cased "bauteil" on some
# sheet).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A1")
    $current = [string]$cell.Value2
    if ($current.ToLower() -eq "bauteil") {
        $cell.Value = "Component"
    }
}

# Keep the originally active sheet selected.
$wb.Worksheets.Item("IfcWall").Activate()
